$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# shifting "Description" and everything below it down by one row.
$ws.Rows.Item(11).EntireRow.Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Bump the published version and date for this ST.r2b update.
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
